$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the LIMITE_FISICO values (column H) for the 5 data rows
$ws.Range("H2:H6").Value = 50

# Remove the CLIENTE_ID column (K) entirely - header + data
$ws.Columns.Item(11).Delete()

# Update selection to match the target state
$ws.Range("L6").Select() | Out-Null
